$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.370.89'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.75%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.773.82'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.41%  '

$ws.Range("E4").Value = '  -0.41%  '

$ws.Range("E5").Value = '  -0.43%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '306.18'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.68%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4231'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.17%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3608'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.32%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07134'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.83%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8375'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.89%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.40'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.97%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.772.95'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.54%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.448'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.56%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.246'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.51%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06895'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.51%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.003'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.62%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '78.92'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.25%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008651'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.75%  '

$ws.Range("E19").Value = '  -0.51%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.90'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.32%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '26.393.49'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.47%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.095'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.86%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.91'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.13%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.993.90'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.03%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.69'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.01%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.796'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -8.17%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.98'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.71%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.059'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.33%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '114.15'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.26%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.840'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +10.64%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08823'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.60%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7263'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.72%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.118'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.61%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.310'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.87%  '

$ws.Range("E35").Value = '  -0.49%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.736'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.47%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.082'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.58%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05105'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.10%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01883'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.70%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.4917'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.52%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1608'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.69%  '

$ws.Range("E42").Value = '  +0.18%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.324'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.27%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.048'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.41%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '104.64'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.03%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.20'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.07%  '

$ws.Range("E47").Value = '  -0.51%  '

$ws.Range("E48").Value = '  +2.29%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06171'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.17%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4439'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.67%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.707'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.06%  '
